$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("sh2")

# --- Sheet1 edits ---
# B1: 10 -> 10.1
$ws1.Range("B1").Value = 10.1

# A3: add formula =3 (value stays 3)
$ws1.Range("A3").Formula = "=3"

# B3: 30 -> string "three" (shared string lookup)
$ws1.Range("B3").Value = "three"

# A4: 4 -> formula =A1+10 (value 11)
$ws1.Range("A4").Formula = "=A1+10"

# A5: 5 -> formula ="this"&A2 (value "thistwo")
$ws1.Range("A5").Formula = '="this"&A2'

# Sheet1 active selection -> B2
$ws1.Range("B2").Select()

# --- sh2 edits ---
# A2: 200 -> string "four"
$ws2.Range("A2").Value = "four"

# sh2 active selection -> A3
$ws2.Range("A3").Select()

# Restore Sheet1 as the active/visible tab (selection on sh2 stays at A3)
$ws1.Activate()
$ws1.Range("B2").Select()

$wb.Save()
